$d = $word.ActiveDocument

# The sentence that needs to be split across five runs.
$target = "Must be a date no later than the Sample Collection Date or the current date."

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target sentence to edit."
}

# Pull the raw WordprocessingML for the document out of the OOXML package
# wrapper that WordOpenXML returns, so we can read back the exact <w:p>
# (paragraph identity/rsid attributes, paragraph formatting) and <w:rPr>
# (run formatting) surrounding the sentence, instead of hard-coding them.
$full = $rng.WordOpenXML
$docPartMarker = '<pkg:part pkg:name="/word/document.xml"'
$partIdx = $full.IndexOf($docPartMarker)
$xmlDataIdx = $full.IndexOf('<pkg:xmlData>', $partIdx)
$docStart = $xmlDataIdx + '<pkg:xmlData>'.Length
$docEnd = $full.IndexOf('</pkg:xmlData>', $docStart)
$docXml = $full.Substring($docStart, $docEnd - $docStart)

$needlePos = $docXml.IndexOf($target)
if ($needlePos -lt 0) {
    throw "Could not locate target sentence in WordOpenXML."
}

# Walk backward from the sentence to the start of its enclosing <w:p>.
$prefix = $docXml.Substring(0, $needlePos)
$pStart = $prefix.LastIndexOf("<w:p ")
if ($pStart -lt 0) {
    $pStart = $prefix.LastIndexOf("<w:p>")
}
$pTagEndRel = $docXml.IndexOf(">", $pStart)
$pOpenTag = $docXml.Substring($pStart, $pTagEndRel - $pStart + 1)

# Grab the <w:pPr>...</w:pPr> block (paragraph formatting) verbatim.
$afterOpenTag = $pTagEndRel + 1
$pPrBlock = ""
if ($docXml.IndexOf("<w:pPr", $afterOpenTag) -eq $afterOpenTag) {
    $pPrEnd = $docXml.IndexOf("</w:pPr>", $afterOpenTag)
    $pPrBlock = $docXml.Substring($afterOpenTag, ($pPrEnd + "</w:pPr>".Length) - $afterOpenTag)
}

# Grab the <w:rPr>...</w:rPr> block (run formatting) of the run holding
# the sentence; it sits between the end of <w:pPr> and the <w:t> element.
$rPrRegionStart = $afterOpenTag + $pPrBlock.Length
$tStart = $docXml.IndexOf("<w:t", $rPrRegionStart)
$runPrefix = $docXml.Substring($rPrRegionStart, $tStart - $rPrRegionStart)
$rPrStart = $runPrefix.IndexOf("<w:rPr>")
$rPr = ""
if ($rPrStart -ge 0) {
    $rPrEndRel = $runPrefix.IndexOf("</w:rPr>") + "</w:rPr>".Length
    $rPr = $runPrefix.Substring($rPrStart, $rPrEndRel - $rPrStart)
}

# Build the replacement paragraph: same <w:p>/<w:pPr>, but the sentence is
# now split across five runs (all keeping the original run formatting),
# reading: "Must be a date no earlier than the Sample Collection Date or
# later than the current date."
$segments = @(
    "Must be a date ",
    "no earlier than ",
    "the Sample Collection Date or ",
    "later than ",
    "the current date."
)

$runsXml = ""
for ($i = 0; $i -lt $segments.Length; $i++) {
    $text = $segments[$i]
    $preserve = ($text -ne $text.Trim())
    $tOpen = "<w:t>"
    if ($preserve) {
        $tOpen = '<w:t xml:space="preserve">'
    }
    $runsXml += "<w:r>" + $rPr + $tOpen + $text + "</w:t></w:r>"
}

$xmlns = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
$pOpenTagWithNs = $pOpenTag.Substring(0, 4) + $xmlns + $pOpenTag.Substring(4)

$xmlFragment = $pOpenTagWithNs + $pPrBlock + $runsXml + "</w:p>"

$rng.InsertXML($xmlFragment) | Out-Null

# Sanity-check: the new sentence should now be found as one continuous
# range of text, and the old sentence should be gone.
$expected = "Must be a date no earlier than the Sample Collection Date or later than the current date."
$check1 = $d.Content.Duplicate
$ok1 = $check1.Find.Execute($expected, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok1) {
    throw "Post-edit verification failed: new sentence not found."
}
$check2 = $d.Content.Duplicate
$ok2 = $check2.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok2) {
    throw "Post-edit verification failed: old sentence still present."
}
